$wb = $excel.ActiveWorkbook

# Update the "runner" sheet's B3 value from "Y" to "N"
$runner = $wb.Worksheets.Item("runner")
$runner.Range("B3").Value = "N"

# Make "runner" the active sheet/tab (was "MyFirstTest")
$runner.Activate()

$wb.Save()
